# Updated Test Case Excel
# Rebuild the test-case table (rows 2-42) to reflect the revised scenarios:
#  - split "add to cart" / "remove from cart" into single- and multiple-product variants
#  - add a new "invalid login" test case
#  - renumber/relayout the remaining rows accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old body (content + formatting) before laying out the new rows so that
# no stale styles/values are left behind at mismatched row numbers.
$ws.Range("A2:L36").Clear()

$ws.Range("A2").Value = "Adding and Deleting Products to Cart"
$ws.Range("B2").Value = 101
$ws.Range("C2").Value = "Verify that a user can successfully add a single product to the cart."
$ws.Range("D2").Value = "test_add_single_product"
$ws.Range("E2").Value = "User must be logged in"
$ws.Range("F2").WrapText = $true
$ws.Range("F2").Value = "Login"
$ws.Range("G2").Value = "Item should be added successfully"
$ws.Range("H2").Value = "Pass"
$ws.Range("I2").Value = "High"
$ws.Range("J2").Value = "Yes"

$ws.Range("F3").Value = "Add item to cart"

$ws.Range("F4").Value = "Verify cart contents"

$ws.Range("B5").Value = 102
$ws.Range("C5").Value = "Verify that a user can successfully add a multiple products to the cart."
$ws.Range("D5").Value = "test_add_multiple_products"
$ws.Range("E5").Value = "User must be logged in"
$ws.Range("F5").WrapText = $true
$ws.Range("F5").Value = "Login"
$ws.Range("G5").Value = "Item should be added successfully"
$ws.Range("H5").Value = "Pass"
$ws.Range("I5").Value = "High"
$ws.Range("J5").Value = "Yes"

$ws.Range("F6").Value = "Select and add multiple items to cart"

$ws.Range("F7").Value = "Verify cart contents"

$ws.Range("B8").Value = 103
$ws.Range("C8").Value = "Verify that a user can remove a product from the cart."
$ws.Range("D8").Value = "test_delete_single_product_from_cart"
$ws.Range("E8").Value = "User must have items in cart"
$ws.Range("F8").WrapText = $true
$ws.Range("F8").Value = "Login"
$ws.Range("G8").Value = "Item should be removed successfully"
$ws.Range("H8").Value = "Pass"
$ws.Range("I8").Value = "High"
$ws.Range("J8").Value = "Yes"

$ws.Range("F9").Value = "Remove item from cart"

$ws.Range("F10").Value = "Verify cart is empty"

$ws.Range("B11").Value = 104
$ws.Range("C11").Value = "Verify that a user can remove multiple products from the cart."
$ws.Range("D11").Value = "test_delete_multiple_products"
$ws.Range("E11").Value = "User must have items in cart"
$ws.Range("F11").WrapText = $true
$ws.Range("F11").Value = "Login"
$ws.Range("G11").Value = "Item should be removed successfully"
$ws.Range("H11").Value = "Pass"
$ws.Range("I11").Value = "High"
$ws.Range("J11").Value = "Yes"

$ws.Range("F12").Value = "Remove all items from cart"

$ws.Range("F13").Value = "Verify cart is empty"

$ws.Range("B14").Value = 105
$ws.Range("C14").Value = "Verify that the cart item details after adding products."
$ws.Range("D14").Value = "test_add_single_product"
$ws.Range("E14").Value = "User must have added items to cart"
$ws.Range("F14").WrapText = $true
$ws.Range("F14").Value = "Login"
$ws.Range("G14").Value = "Total Cart Items details must display correctly"
$ws.Range("H14").Value = "Pass"
$ws.Range("I14").Value = "Medium"
$ws.Range("J14").Value = "Yes"

$ws.Range("F15").Value = "Add multiple items"

$ws.Range("F16").Value = "Verify total price"

$ws.Range("B17").Value = 106
$ws.Range("C17").Value = "Verify that the cart displays a message when empty after removing all products."
$ws.Range("D17").Value = "test_delete_multiple_products"
$ws.Range("E17").Value = "User must have added items to cart"
$ws.Range("F17").WrapText = $true
$ws.Range("F17").Value = "Login"
$ws.Range("G17").Value = "Empty cart message should be displayed"
$ws.Range("H17").Value = "Pass"
$ws.Range("I17").Value = "Medium"
$ws.Range("J17").Value = "Yes"

$ws.Range("F18").Value = "Remove all items"

$ws.Range("F19").Value = "Check empty cart message"

$ws.Range("A20").Value = "Registration and Login"
$ws.Range("B20").Value = 201
$ws.Range("C20").Value = "Verify that a user can successfully register with valid credentials."
$ws.Range("D20").Value = "test_user_registration"
$ws.Range("E20").Value = "None"
$ws.Range("F20").WrapText = $true
$ws.Range("F20").Value = "Open registration page"
$ws.Range("G20").Value = "User should be registered successfully"
$ws.Range("H20").Value = "Pass"
$ws.Range("I20").Value = "High"
$ws.Range("J20").Value = "Yes"

$ws.Range("F21").Value = "Enter valid details"

$ws.Range("F22").Value = "Submit and verify"

$ws.Range("B23").Value = 202
$ws.Range("C23").Value = "Verify that a registered user can log in with valid credentials."
$ws.Range("D23").Value = "test_login"
$ws.Range("E23").Value = "User must be registered"
$ws.Range("F23").WrapText = $true
$ws.Range("F23").Value = "Open login page"
$ws.Range("G23").Value = "User should log in successfully"
$ws.Range("H23").Value = "Pass"
$ws.Range("I23").Value = "High"
$ws.Range("J23").Value = "Yes"

$ws.Range("F24").Value = "Enter valid credentials"

$ws.Range("F25").Value = "Verify login success"

$ws.Range("B26").Value = 203
$ws.Range("C26").Value = "Verify that login fails with incorrect credentials."
$ws.Range("D26").Value = "test_invalid_login"
$ws.Range("E26").Value = "User must be registered"
$ws.Range("F26").WrapText = $true
$ws.Range("F26").Value = "Open login page"
$ws.Range("G26").Value = "Error message should appear"
$ws.Range("H26").Value = "Pass"
$ws.Range("I26").Value = "Medium"
$ws.Range("J26").Value = "Yes"

$ws.Range("F27").Value = "Enter incorrect credentials"

$ws.Range("F28").Value = "Verify error message"

$ws.Range("B29").Value = 204
$ws.Range("C29").Value = "Verify that a user can log out successfully."
$ws.Range("D29").Value = "test_login"
$ws.Range("E29").Value = "User must be logged in"
$ws.Range("F29").WrapText = $true
$ws.Range("F29").Value = "Login"
$ws.Range("G29").Value = "User should be logged out"
$ws.Range("H29").Value = "Pass"
$ws.Range("I29").Value = "Medium"
$ws.Range("J29").Value = "Yes"

$ws.Range("F30").Value = "Click logout"

$ws.Range("F31").Value = "Verify user is logged out"

$ws.Range("A32").Value = "Add Address to Your Profile and Verify"
$ws.Range("B32").Value = 301
$ws.Range("C32").Value = "Verify that a user can successfully add a billing address."
$ws.Range("D32").Value = "test_add_address"
$ws.Range("E32").Value = "User must be logged in"
$ws.Range("F32").WrapText = $true
$ws.Range("F32").Value = "Login"
$ws.Range("G32").Value = "Billing address should be saved"
$ws.Range("H32").Value = "Pass"
$ws.Range("I32").Value = "High"
$ws.Range("J32").Value = "Yes"

$ws.Range("F33").Value = "Navigate to address section"

$ws.Range("F34").Value = "Add billing address"

$ws.Range("F35").Value = "Save address"

$ws.Range("B36").Value = 302
$ws.Range("C36").Value = "Verify that the saved billing address appears correctly in the profile."
$ws.Range("D36").Value = "test_saved_address"
$ws.Range("E36").Value = "User must have added an address"
$ws.Range("F36").WrapText = $true
$ws.Range("F36").Value = "Login"
$ws.Range("G36").Value = "Saved address should match to expected address"
$ws.Range("H36").Value = "Pass"
$ws.Range("I36").Value = "Medium"
$ws.Range("J36").Value = "Yes"

$ws.Range("F37").Value = "Navigate to profile"

$ws.Range("F38").Value = "Verify saved address"

$ws.Range("B39").Value = 303
$ws.Range("C39").Value = "Verify that the system prevents saving an address with missing required fields (e.g., first name, city)."
$ws.Range("D39").Value = "NA"
$ws.Range("E39").Value = "User must be logged in"
$ws.Range("F39").WrapText = $true
$ws.Range("F39").Value = "Login"
$ws.Range("G39").Value = "System should prevent saving incomplete address"
$ws.Range("H39").Value = "Not Run"
$ws.Range("I39").Value = "High"
$ws.Range("J39").Value = "No"

$ws.Range("F40").Value = "Navigate to address section"

$ws.Range("F41").Value = "Enter incomplete address"

$ws.Range("F42").Value = "Verify error message"

# Column F ("Test Steps") needs to be a bit wider to fit the new, longer step text.
$ws.Columns.Item(6).ColumnWidth = 30.6666667

# Restore the active selection and page orientation.
$ws.Range("C9").Select()
$ws.PageSetup.Orientation = 1

